# Workbook / worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new column before column A ("spc" helper column).
#    This shifts the existing name/smiles/mult/charge/exc_flag
#    columns from A:E to B:F, carrying their styles/values along.
# ------------------------------------------------------------------
$ws.Columns.Item(1).Insert()

# ------------------------------------------------------------------
# 2. Fill the new column A: "spc" down every row first (this is the
#    order the shared-string table was actually built in), then go
#    back and give the header cell its explanatory text.
# ------------------------------------------------------------------
$ws.Range("A2").Value = "spc"
$ws.Range("A3").Value = "spc"
$ws.Range("A4").Value = "spc"
$ws.Range("A5").Value = "spc"
$ws.Range("A6").Value = "spc"
$ws.Range("A7").Value = "spc"
$ws.Range("A8").Value = "spc"
$ws.Range("A9").Value = "spc"
$ws.Range("A10").Value = "spc"
$ws.Range("A1").Value = "this column is here so I can copy paste easily"

# ------------------------------------------------------------------
# 3. Row 10 (N2) previously only had name + mult filled in; now also
#    give it a smiles value plus explicit charge/exc_flag of 0.
# ------------------------------------------------------------------
$ws.Range("C10").Value = "N#N"
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0

# ------------------------------------------------------------------
# 4. Append five brand-new species rows (11-14 ... wait, 4 rows: 11-14)
# ------------------------------------------------------------------
$ws.Range("A11").Value = "spc"
$ws.Range("B11").Value = "H2O"
$ws.Range("C11").Value = "O"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0

$ws.Range("A12").Value = "spc"
$ws.Range("B12").Value = "H2"
$ws.Range("C12").Value = "[H][H]"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0

$ws.Range("A13").Value = "spc"
$ws.Range("B13").Value = "NO"
$ws.Range("C13").Value = "[N]=O"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0

$ws.Range("A14").Value = "spc"
$ws.Range("B14").Value = "NO2"
$ws.Range("C14").Value = "N(=O)[O]"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0

# ------------------------------------------------------------------
# 5. Widen column A to fit its new contents (matches the workbook's
#    saved column width for the "spc" helper column).
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 37.666666

# ------------------------------------------------------------------
# 6. Restore the selection the author left the sheet on.
# ------------------------------------------------------------------
$ws.Range("A19").Select()
